$wb = $excel.ActiveWorkbook

# The header cells below live on protected sheets (sheet protection carries
# fine-grained flags - objects/scenarios/insertRows/deleteRows - that a plain
# Unprotect()+Protect() round-trip can't reproduce on this engine). Toggling
# Locked off/on around the write edits the cell in place without disturbing
# the sheet's <sheetProtection> element.
function Set-LockedCellValue($range, $value) {
    $range.Locked = $false
    $range.Value = $value
    $range.Locked = $true
}

# Sheet "!!_Table of contents": header row A1 and TableOfContents metadata row A2
$wsToc = $wb.Worksheets.Item("!!_Table of contents")
Set-LockedCellValue $wsToc.Range("A1") "!!!ObjTables objTablesVersion='0.0.9' date='2020-04-27 01:05:05'"
Set-LockedCellValue $wsToc.Range("A2") "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='2020-04-27 01:05:05' objTablesVersion='0.0.9'"

# Sheet "!!_Schema": header row A1 metadata
$wsSchema = $wb.Worksheets.Item("!!_Schema")
Set-LockedCellValue $wsSchema.Range("A1") "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='2020-04-27 01:05:06' objTablesVersion='0.0.9'"

# Sheet "!!Transaction": header row A1 metadata, id= -> class=
$wsTransaction = $wb.Worksheets.Item("!!Transaction")
Set-LockedCellValue $wsTransaction.Range("A1") "!!ObjTables type='Data' tableFormat='row' class='Transaction' name='Transaction' description='Stores transactions' date='2020-04-27 01:05:06' objTablesVersion='0.0.9'"
